# New Test Case added: TCR009 and TCR010 rows on the "Register" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Register")

# ---------------------------------------------------------------------
# Row 10 - TCR009 "Verify whether the password fields ... Complexity Standards"
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "TCR009"
$ws.Range("B10").Value = "TS001"
$ws.Range("C10").Value = "Register"
$ws.Range("D10").Value = "Verify whether the password fields in the register account page are following Password Complexity Standards"
$ws.Range("E10").Value = "1. Application (https://tutorialsninja.com/demo) is opened"
$ws.Range("F10").Value = "1. Click on 'My Account' dropdown menu`n2. Click on 'Register' option `n3. Enter new account details into all the fields (First Name, Last Name, E-Mail, Telephone, Password, Password Confirm, Newsletter and  Privacy Policy fields)`n4. Check entering simple passwords (Not following Password Complexity Standars' i.e. Size of password as 8, password should contain atleast one number, symbol, lower case letter and upper case letters) - <Refer Test Data>`n5. Click on 'Continue' button (ER-1)"
$ws.Range("G10").Value = "Try all below passwords:`n1) 12345`n2) abcde`n"
$ws.Range("H10").Value = "1. Warning message should be displayed for following Password Complexity Standards"

# ---------------------------------------------------------------------
# Row 11 - TCR010 "Verify registering the account without selecting 'Privacy Policy'"
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "TCR010"
$ws.Range("B11").Value = "TS001"
$ws.Range("C11").Value = "Register"
$ws.Range("D11").Value = "Verify registering the account without selecting the 'Privacy Policy' checkbox option"
$ws.Range("E11").Value = "1. Application (https://tutorialsninja.com/demo) is opened"
$ws.Range("F11").Value = "1. Click on 'My Account' dropdown menu`n2. Click on 'Register' option `n3. Enter new account details into all the fields (First Name, Last Name, E-Mail, Telephone, Password, Password Confirm and Newsletter fields)`n4. Don't select the 'Privacy Policy' checkbox option`n4. Click on 'Continue' button (ER-1)"
$ws.Range("G11").Value = "Not Applicable"
$ws.Range("H11").Value = "1. Warning message - 'Warning: You must agree to the Privacy Policy!' should be displayed"

# ---------------------------------------------------------------------
# Formatting: copy the look of the existing populated rows (8/9) onto the
# two new rows, cell by cell, so borders/fonts/fills/alignment match the
# rest of the table. (Row 9 already has the exact formatting rows 10/11
# need for most columns; G11 needs the "centered" look used elsewhere for
# "Not Applicable" data cells, e.g. C12.)
# ---------------------------------------------------------------------
function Copy-Format($srcAddr, $dstAddr) {
    $src = $ws.Range($srcAddr)
    $dst = $ws.Range($dstAddr)
    $src.Copy()
    $dst.PasteSpecial(-4122) # xlPasteFormats
}

Copy-Format "A9" "A10"
Copy-Format "B9" "B10"
Copy-Format "C9" "C10"
Copy-Format "E9" "E10"
Copy-Format "F9" "F10"
Copy-Format "G9" "G10"
Copy-Format "H9" "H10"

Copy-Format "A9" "A11"
Copy-Format "B9" "B11"
Copy-Format "C9" "C11"
Copy-Format "E9" "E11"
Copy-Format "F9" "F11"
Copy-Format "H9" "H11"
Copy-Format "C12" "G11"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row heights for the two new rows
# ---------------------------------------------------------------------
$ws.Rows.Item(10).RowHeight = 218.4
$ws.Rows.Item(11).RowHeight = 140.4

# ---------------------------------------------------------------------
# Sheet view: scroll down a bit and move the selection to F11
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F11").Select()
